# Add a team "record" (Wins / Losses / Ties) to the roster sheet.
# New columns AD, AE, AF are appended after the existing AC column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row (data currently spans A1:AC67).
$lastRow = $ws.UsedRange.Rows.Count

# Column numbers for AD, AE, AF.
$colWins   = 30
$colLosses = 31
$colTies   = 32

# --- Header row (row 1) ---------------------------------------------------
$ws.Cells.Item(1, $colWins).Value2   = "Wins"
$ws.Cells.Item(1, $colLosses).Value2 = "Losses"
$ws.Cells.Item(1, $colTies).Value2   = "Ties"

# Match the existing header formatting (bold, centered, bordered) by
# copying the style from the neighboring header cell (AC1).
$ws.Cells.Item(1, 29).Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (2..last row) ----------------------------------------------
# Every player row gets the same team record: 82 wins, 80 losses, 0 ties.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colWins).Value2   = 82
    $ws.Cells.Item($r, $colLosses).Value2 = 80
    $ws.Cells.Item($r, $colTies).Value2   = 0
}
